$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at position 538, pushing existing rows 538-615 down to 542-619
$ws.Rows("538:541").Insert()

# Common/constant column values for this product across all rows in this block
$marketId = 9
$market = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$catId = 100112037
$cat = "Cebollín"
$variedad = "Sin especificar"
$unidad = "`$/paquete 36 unidades"
$origen = "Región Metropolitana"
$kgUnidades = 36
$clasif = "Hortaliza"

function Set-DataRow($r, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $precioKg) {
    $ws.Cells.Item($r, 1).Value = $marketId
    $ws.Cells.Item($r, 2).Value = $market
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $catId
    $ws.Cells.Item($r, 7).Value = $cat
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $pmin
    $ws.Cells.Item($r, 12).Value = $pmax
    $ws.Cells.Item($r, 13).Value = $pprom
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $precioKg
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasif
}

Set-DataRow 538 44474 "Extra"   106 2400 2600 2500 69
Set-DataRow 539 44474 "Primera" 250 1900 2200 2050 57
Set-DataRow 540 44474 "Segunda" 160 1600 1800 1700 47
Set-DataRow 541 44474 "Tercera" 97  1200 1200 1200 33
